$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for the removed worker entries (rows 17-22, OSWALDO ENRIQUE CARBALLO DE VOZ)
$ws.Range("B17:J22").EntireRow.Delete() | Out-Null

